# Refresh cryptos list values (prices + 1h volume %) per the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.605.80'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '3.084.55'
$ws.Range('E3').Value = '  -0.31%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = '''543.56'
$ws.Range('E5').Value = '  -1.28%  '

$ws.Range('D6').Value = '''140.22'
$ws.Range('E6').Value = '  +2.04%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').Value = '3.079.56'
$ws.Range('E8').Value = '  -0.30%  '

$ws.Range('D9').Value = '''0.501'
$ws.Range('E9').Value = '  +0.81%  '

$ws.Range('E10').Value = '  -1.11%  '

$ws.Range('D11').Value = '''6.41'
$ws.Range('E11').Value = '  +1.07%  '

$ws.Range('D12').Value = '''0.457'
$ws.Range('E12').Value = '  -2.65%  '

$ws.Range('D13').Value = '''0.0000224'
$ws.Range('E13').Value = '  +2.92%  '

$ws.Range('D14').Value = '''34.92'
$ws.Range('E14').Value = '  -1.51%  '

$ws.Range('D15').Value = '3.582.39'
$ws.Range('E15').Value = '  -0.32%  '

$ws.Range('D16').Value = '63.563.46'
$ws.Range('E16').Value = '  +0.27%  '

$ws.Range('E17').Value = '  +0.98%  '

$ws.Range('D18').Value = '3.080.83'
$ws.Range('E18').Value = '  -0.48%  '

$ws.Range('D19').Value = '''6.64'
$ws.Range('E19').Value = '  -1.65%  '

$ws.Range('D20').Value = '''473.66'
$ws.Range('E20').Value = '  -3.58%  '

$ws.Range('D21').Value = '''13.45'
$ws.Range('E21').Value = '  -1.66%  '

$ws.Range('D22').Value = '''0.700'
$ws.Range('E22').Value = '  -2.26%  '

$ws.Range('D23').Value = '''7.10'
$ws.Range('E23').Value = '  -2.07%  '

$ws.Range('D24').Value = '''78.81'
$ws.Range('E24').Value = '  -0.25%  '

$ws.Range('D25').Value = '''12.27'
$ws.Range('E25').Value = '  -1.17%  '

$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('D27').Value = '''2.71'
$ws.Range('E27').Value = '  -1.71%  '

$ws.Range('D28').Value = '''8.02'
$ws.Range('E28').Value = '  -5.29%  '

$ws.Range('D29').Value = '''0.997'
$ws.Range('E29').Value = '  -0.24%  '

$ws.Range('D30').Value = '''26.24'
$ws.Range('E30').Value = '  -1.38%  '

$ws.Range('D31').Value = '''1.91'
$ws.Range('E31').Value = '  -4.07%  '

$ws.Range('E32').Value = '  +1.60%  '

$ws.Range('D33').Value = '''57.81'
$ws.Range('E33').Value = '  -1.81%  '

$ws.Range('D34').Value = '''2.34'
$ws.Range('E34').Value = '  -7.04%  '

$ws.Range('D35').Value = '''5.46'
$ws.Range('E35').Value = '  +5.85%  '

$ws.Range('D36').Value = '''493.71'
$ws.Range('E36').Value = '  -5.16%  '

$ws.Range('D37').Value = '''6.01'
$ws.Range('E37').Value = '  -0.04%  '

$ws.Range('D38').Value = '3.246.35'
$ws.Range('E38').Value = '  +3.29%  '

$ws.Range('D39').Value = '''0.0403'
$ws.Range('E39').Value = '  -1.04%  '

$ws.Range('D40').Value = '''0.0799'
$ws.Range('E40').Value = '  -0.88%  '

$ws.Range('E41').Value = '  -0.54%  '

$ws.Range('D42').Value = '''8.12'
$ws.Range('E42').Value = '  -0.80%  '

$ws.Range('D43').Value = '''2.66'
$ws.Range('E43').Value = '  -0.80%  '

$ws.Range('D44').Value = '''0.255'
$ws.Range('E44').Value = '  -1.28%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '''25.42'
$ws.Range('E46').Value = '  +1.05%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '''123.92'
$ws.Range('E47').Value = '  +1.91%  '

$ws.Range('D48').Value = '''2.05'
$ws.Range('E48').Value = '  -1.75%  '

$ws.Range('D49').Value = '0.0₃0533'
$ws.Range('E49').Value = '  +5.51%  '

$ws.Range('E50').Value = '  +0.32%  '

$ws.Range('E51').Value = '  +4.46%  '
